$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'2019-12-3"
$ws.Range("E2").Value = "'115585000"
$ws.Range("D3").Value = "'2015-12-3"
$ws.Range("E3").Value = "'118779000"
$ws.Range("D4").Value = "'2018-12-3"
$ws.Range("E4").Value = "'113923000"
$ws.Range("D5").Value = "'2016-12-3"
$ws.Range("E5").Value = "'105590000"
$ws.Range("D7").Value = "'2015-12-3"
$ws.Range("E7").Value = "'0.067"
$ws.Range("D8").Value = "'2017-12-3"
$ws.Range("E8").Value = "'0.028"
$ws.Range("D9").Value = "'2016-12-3"
$ws.Range("D10").Value = "'2015-12-3"
$ws.Range("E10").Value = "'0.032"
$ws.Range("D11").Value = "'2017-12-3"
$ws.Range("E11").Value = "'0.032"
$ws.Range("D12").Value = "'2019-12-3"
$ws.Range("D13").Value = "'2018-12-3"
$ws.Range("E13").Value = "'0.028"
$ws.Range("D15").Value = "'2016-12-3"
$ws.Range("E15").Value = "'3701000"
$ws.Range("D16").Value = "'2015-12-3"
$ws.Range("E16").Value = "'8822000"
$ws.Range("D17").Value = "'2019-12-3"
$ws.Range("E17").Value = "'3136000"
$ws.Range("D18").Value = "'2018-12-3"
$ws.Range("E18").Value = "'6140000"
